$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2078.7693
$ws.Range("I62").Value = 1898.4736
$ws.Range("J62").Value = 2568.1428
$ws.Range("K62").Value = 1898.4736
$ws.Range("L62").Value = 2568.1428
$ws.Range("M62").Value = -1274.4736
$ws.Range("N62").Value = -3816.1428
$ws.Range("H65").Value = 2078.7693
$ws.Range("I65").Value = 1898.4736
$ws.Range("J65").Value = 2568.1428
$ws.Range("K65").Value = 9492.368
$ws.Range("L65").Value = 12840.714
$ws.Range("M65").Value = -6372.368
$ws.Range("N65").Value = -19080.714
$ws.Range("H98").Value = 1554
$ws.Range("I98").Value = 1488.8889
$ws.Range("J98").Value = 1749.3334
$ws.Range("K98").Value = 1488.8889
$ws.Range("L98").Value = 1749.3334
$ws.Range("M98").Value = 9.111100000000079
$ws.Range("N98").Value = -4745.3334
$ws.Range("H106").Value = 3827.1365
$ws.Range("I106").Value = 4212.3335
$ws.Range("J106").Value = 3001.7144
$ws.Range("K106").Value = 4212.3335
$ws.Range("L106").Value = 3001.7144
$ws.Range("M106").Value = -3581.3335
$ws.Range("N106").Value = -4263.7144
$ws.Range("H112").Value = 2037.04
$ws.Range("J112").Value = 2101.0833
$ws.Range("L112").Value = 6303.249899999999
$ws.Range("N112").Value = -8519.249899999999
$ws.Range("H113").Value = 1923.125
$ws.Range("I113").Value = 1884.1666
$ws.Range("J113").Value = 2040
$ws.Range("K113").Value = 1884.1666
$ws.Range("L113").Value = 2040
$ws.Range("M113").Value = 1369.8334
$ws.Range("N113").Value = -8548
$ws.Range("H116").Value = 1478.4667
$ws.Range("I116").Value = 1482.3636
$ws.Range("J116").Value = 1467.75
$ws.Range("K116").Value = 1482.3636
$ws.Range("L116").Value = 1467.75
$ws.Range("M116").Value = 1959.6364
$ws.Range("N116").Value = -8351.75
$ws.Range("H121").Value = 725.2
$ws.Range("I121").Value = 645
$ws.Range("J121").Value = 732.1739
$ws.Range("K121").Value = 1935
$ws.Range("L121").Value = 2196.5217
$ws.Range("M121").Value = -188
$ws.Range("N121").Value = -5690.5217
$ws.Range("H122").Value = 1554
$ws.Range("I122").Value = 1488.8889
$ws.Range("J122").Value = 1749.3334
$ws.Range("K122").Value = 4466.6667
$ws.Range("L122").Value = 5248.0002
$ws.Range("M122").Value = -2016.6667
$ws.Range("N122").Value = -10148.0002
$ws.Range("H129").Value = 862.9231
$ws.Range("I129").Value = 419.72726
$ws.Range("J129").Value = 953.2037
$ws.Range("K129").Value = 1259.18178
$ws.Range("L129").Value = 2859.6111
$ws.Range("M129").Value = 3740.81822
$ws.Range("N129").Value = -12859.6111
$ws.Range("H132").Value = 1443637
$ws.Range("I132").Value = 2567.4062
$ws.Range("J132").Value = 24500750
$ws.Range("K132").Value = 7702.2186
$ws.Range("L132").Value = 73502250
$ws.Range("M132").Value = -5172.2186
$ws.Range("N132").Value = -73507310
$ws.Range("H137").Value = 3032247
$ws.Range("I137").Value = 5557239.5
$ws.Range("J137").Value = 2256
$ws.Range("K137").Value = 16671718.5
$ws.Range("L137").Value = 6768
$ws.Range("M137").Value = -16669168.5
$ws.Range("N137").Value = -11868
$ws.Range("H138").Value = 2419896.8
$ws.Range("I138").Value = 2254.5715
$ws.Range("J138").Value = 3477615.2
$ws.Range("K138").Value = 6763.7145
$ws.Range("L138").Value = 10432845.6
$ws.Range("M138").Value = -1623.7145
$ws.Range("N138").Value = -10443125.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3636801.5
$ws.Range("I32").Value = 3889739.5
$ws.Range("K32").Value = 3889739.5
$ws.Range("M32").Value = -3889452.5
$ws.Range("H74").Value = 7200933.5
$ws.Range("I74").Value = 12550701
$ws.Range("J74").Value = 67910
$ws.Range("K74").Value = 12550701
$ws.Range("L74").Value = 67910
$ws.Range("M74").Value = -12549827
$ws.Range("N74").Value = -69658
$ws.Range("H77").Value = 7200933.5
$ws.Range("I77").Value = 12550701
$ws.Range("J77").Value = 67910
$ws.Range("K77").Value = 62753505
$ws.Range("L77").Value = 339550
$ws.Range("M77").Value = -62749137
$ws.Range("N77").Value = -348286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 71431060
$ws.Range("I105").Value = 83335480
$ws.Range("K105").Value = 83335480
$ws.Range("M105").Value = -83333733

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1904
$ws.Range("I105").Value = 2004.4445
$ws.Range("K105").Value = 2004.4445
$ws.Range("M105").Value = -257.4445000000001
$ws.Range("H132").Value = 44405.457
$ws.Range("I132").Value = 2690.1765
$ws.Range("J132").Value = 145714
$ws.Range("K132").Value = 8070.529500000001
$ws.Range("L132").Value = 437142
$ws.Range("M132").Value = -5540.529500000001
$ws.Range("N132").Value = -442202

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 19402.17
$ws.Range("I5").Value = 36115.895
$ws.Range("J5").Value = 682.8
$ws.Range("K5").Value = 108347.685
$ws.Range("L5").Value = 2048.4
$ws.Range("M5").Value = -108235.685
$ws.Range("N5").Value = -2272.4
$ws.Range("H68").Value = 1228.0615
$ws.Range("I68").Value = 638.2727
$ws.Range("J68").Value = 1529.814
$ws.Range("K68").Value = 1914.8181
$ws.Range("L68").Value = 4589.442
$ws.Range("M68").Value = -1103.8181
$ws.Range("N68").Value = -6211.442
$ws.Range("H71").Value = 1228.0615
$ws.Range("I71").Value = 638.2727
$ws.Range("J71").Value = 1529.814
$ws.Range("K71").Value = 5744.454299999999
$ws.Range("L71").Value = 13768.326
$ws.Range("M71").Value = -1688.454299999999
$ws.Range("N71").Value = -21880.326
$ws.Range("H113").Value = 572.7560999999999
$ws.Range("I113").Value = 524.6667
$ws.Range("J113").Value = 600.5
$ws.Range("K113").Value = 1574.0001
$ws.Range("L113").Value = 1801.5
$ws.Range("M113").Value = 595.9999
$ws.Range("N113").Value = -6141.5
$ws.Range("H122").Value = 522.7619
$ws.Range("I122").Value = 280.88235
$ws.Range("J122").Value = 1550.75
$ws.Range("K122").Value = 2527.94115
$ws.Range("L122").Value = 13956.75
$ws.Range("M122").Value = -77.94114999999965
$ws.Range("N122").Value = -18856.75
$ws.Range("H135").Value = 19402.17
$ws.Range("I135").Value = 36115.895
$ws.Range("J135").Value = 682.8
$ws.Range("K135").Value = 325043.055
$ws.Range("L135").Value = 6145.2
$ws.Range("M135").Value = -322508.055
$ws.Range("N135").Value = -11215.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 6214
$ws.Range("I99").Value = 4527.9287
$ws.Range("J99").Value = 29819
$ws.Range("K99").Value = 4527.9287
$ws.Range("L99").Value = 29819
$ws.Range("M99").Value = -2281.9287
$ws.Range("N99").Value = -34311
$ws.Range("H122").Value = 2905.6667
$ws.Range("I122").Value = 2646.7273
$ws.Range("K122").Value = 7940.1819
$ws.Range("M122").Value = -5490.1819
$ws.Range("H132").Value = 41397.98
$ws.Range("I132").Value = 30442.371
$ws.Range("J132").Value = 65363.375
$ws.Range("K132").Value = 91327.113
$ws.Range("L132").Value = 196090.125
$ws.Range("M132").Value = -88797.113
$ws.Range("N132").Value = -201150.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 60405.883
$ws.Range("J100").Value = 1642.8572
$ws.Range("L100").Value = 1642.8572
$ws.Range("N100").Value = -2724.8572
$ws.Range("H122").Value = 3434.7646
$ws.Range("I122").Value = 2924.3333
$ws.Range("J122").Value = 4659.8
$ws.Range("K122").Value = 8772.999899999999
$ws.Range("L122").Value = 13979.4
$ws.Range("M122").Value = -6322.999899999999
$ws.Range("N122").Value = -18879.4
$ws.Range("H132").Value = 26031.166
$ws.Range("I132").Value = 1530.4783
$ws.Range("J132").Value = 55689.895
$ws.Range("K132").Value = 4591.4349
$ws.Range("L132").Value = 167069.685
$ws.Range("M132").Value = -2061.4349
$ws.Range("N132").Value = -172129.685
$ws.Range("H136").Value = 57638.902
$ws.Range("I136").Value = 26951.6
$ws.Range("J136").Value = 159929.92
$ws.Range("K136").Value = 80854.79999999999
$ws.Range("L136").Value = 479789.76
$ws.Range("M136").Value = -78304.79999999999
$ws.Range("N136").Value = -484889.76

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 57497.832
$ws.Range("I132").Value = 36359.1
$ws.Range("J132").Value = 145072.58
$ws.Range("K132").Value = 109077.3
$ws.Range("L132").Value = 435217.74
$ws.Range("M132").Value = -106547.3
$ws.Range("N132").Value = -440277.74
$ws.Range("H136").Value = 50777.285
$ws.Range("I136").Value = 31895.205
$ws.Range("K136").Value = 95685.61500000001
$ws.Range("M136").Value = -93135.61500000001
